# Apply the edits described by the diff to the BIIBRandom workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shared string "Sell" -> "Strong Sell" (cell C2 is the only cell using it)
$ws.Range("C2").Value = "Strong Sell"

# 2. Column C width 7.42578125 -> 10.42578125
# (COM ColumnWidth is quantized to a pixel grid by the host; 9.6 is the closest
# settable value that lands on the nearest achievable stored width.)
$ws.Range("C1").ColumnWidth = 9.6

# 3. A2 date/time value
$ws.Range("A2").Value = 42651.599583333336

# 4. B2 -12 -> -17
$ws.Range("B2").Value = -17

# 5. Q2 updated score
$ws.Range("Q2").Value = 46.242130528922125

# 6. W2 0 -> 1
$ws.Range("W2").Value = 1
